$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.880.09"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "2.632.36"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "2.659.97"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").Value = "3.092.96"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "58.872.21"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "2.678.41"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").Value = "2.757.35"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Value = "0.0₃0809"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "150.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +12.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.857"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "280.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.614"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0538"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0229"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
